$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 306, pushing existing rows 306-326 down to 307-327
$ws.Rows.Item(306).Insert()

# Populate the newly inserted row 306 with the new record's data
$ws.Range("A306").Value = 10
$ws.Range("B306").Value = "Vega Modelo de Temuco"
$ws.Range("C306").Value = "La Araucanía"
$ws.Range("D306").Value = 45021
$ws.Range("E306").Value = 9
$ws.Range("F306").Value = 100112043
$ws.Range("G306").Value = "Pepino dulce"
$ws.Range("H306").Value = "Cultivar IV Región"
$ws.Range("I306").Value = "Primera"
$ws.Range("J306").Value = 35
$ws.Range("K306").Value = 17000
$ws.Range("L306").Value = 17000
$ws.Range("M306").Value = 17000
$ws.Range("N306").Value = "$/bandeja 18 kilos"
$ws.Range("O306").Value = "Provincia de Limarí"
$ws.Range("P306").Value = 944
$ws.Range("Q306").Value = 18
$ws.Range("R306").Value = "Hortaliza"
